# Apply the edits described by the diff:
#  - Metadata sheet: update URL, Version, Date, Publisher values
#  - Elements sheet: clear the Constraint(s) value for the base "Extension" row (row 2)

$wb = $excel.ActiveWorkbook

# --- Metadata sheet -------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/short-term-care-waiting-period"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

# --- Elements sheet --------------------------------------------------------
$elements = $wb.Worksheets.Item("Elements")

# Row 2 is the base "Extension" element; column AI is "Constraint(s)".
$elements.Range("AI2").Value = ""
